# Deliverable 1 for Project 2 - content update
$d = $word.ActiveDocument

# 1. Remove the stray "_GoBack" bookmark left over from the previous save
#    (Word drops this automatically whenever the file is saved from a
#    session that never used it to navigate back to).
try {
    $gb = $d.Bookmarks.Item("_GoBack")
    if ($gb -ne $null) { $gb.Delete() }
} catch {
    # no-op if the bookmark isn't present
}

# 2. "Name:" -> "Name: Eric Butler"
$pName = $d.Paragraphs.Item(4)
$rName = $pName.Range
$rName.End = $rName.End - 1
if ($rName.Text -eq "Name:") {
    $rName.InsertAfter(" Eric Butler")
}

# 3. "Student Number:" -> "Student Number: 20094078"
$pStudent = $d.Paragraphs.Item(5)
$rStudent = $pStudent.Range
$rStudent.End = $rStudent.End - 1
if ($rStudent.Text -eq "Student Number:") {
    $rStudent.InsertAfter(" 20094078")
}

# 4. "Titles:" -> "Commercial Title: Incorrect Direction Assist (IDA)"
$pTitles = $d.Paragraphs.Item(6)
$rTitles = $pTitles.Range
$rTitles.End = $rTitles.End - 1
[void]$rTitles.Find.Execute("Titles:", $true, $false, $false, $false, $false, $true, 1, $false, `
    "Commercial Title: Incorrect Direction Assist (IDA)", 2)

# 5. Insert three new paragraphs after the title paragraph, before "Profile photo:"
$pTitles = $d.Paragraphs.Item(6)
$pTitles.Range.InsertParagraphAfter()

$pAcademic = $d.Paragraphs.Item(7)
$pAcademic.Range.InsertAfter("Academic Title: Incorrect Direction Detection and Assistance, An ADAS concept")
$pAcademic.Range.InsertParagraphAfter()

$pLanding = $d.Paragraphs.Item(8)
$pLanding.Range.InsertAfter("Landing page URL: https://ericdebuitleir.github.io/FYP-WebPage/")
$pLanding.Range.InsertParagraphAfter()

$pDesc = $d.Paragraphs.Item(9)
$pDesc.Range.InsertAfter("Description: Incorrect Direction Assist is a simulation project based on the CANoe software architecture that will detect a vehicle going the incorrect way through a route and deploy aids to correct this based upon current ADAS systems. The ADAS systems simulated here include EBA (Emergency brake assist), Blinking LED from a Dashboard (Simulate RADAR sensors), LA (Lane Assist)")
